$wb = $excel.ActiveWorkbook

# --- "parts_1" sheet: update the Ponoko "Case" price from 57.92 to 60 ---
$ws1 = $wb.Worksheets.Item("parts_1")
$ws1.Range("C2").Value = 60

# --- "parts_10" sheet: update the Ponoko price-alternative note on the Wood case row ---
$ws2 = $wb.Worksheets.Item("parts_10")
$ws2.Range("G4").Value = "Or order 10+ from Ponoko for @`$48"

# --- update the selected cell on each sheet to match the saved view state ---
$ws1.Activate()
$ws1.Range("C3").Select() | Out-Null

$ws2.Activate()
$ws2.Range("A5").Select() | Out-Null
